$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H69").Value = 10899.8
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 11499.667
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 34499.001
$ws.Range("M69").Value = -29126
$ws.Range("N69").Value = -36247.001
$ws.Range("H70").Value = 682781.1
$ws.Range("J70").Value = 2956.0908
$ws.Range("L70").Value = 8868.2724
$ws.Range("N70").Value = -9408.2724
$ws.Range("H72").Value = 10899.8
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 11499.667
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 103497.003
$ws.Range("M72").Value = -85632
$ws.Range("N72").Value = -112233.003
$ws.Range("H73").Value = 682781.1
$ws.Range("J73").Value = 2956.0908
$ws.Range("L73").Value = 8868.2724
$ws.Range("N73").Value = -10740.2724
$ws.Range("H80").Value = 635661.3
$ws.Range("J80").Value = 1671.6666
$ws.Range("L80").Value = 5014.9998
$ws.Range("N80").Value = -7010.9998
$ws.Range("H83").Value = 635661.3
$ws.Range("J83").Value = 1671.6666
$ws.Range("L83").Value = 15044.9994
$ws.Range("N83").Value = -25028.9994
$ws.Range("H127").Value = 1170.7693
$ws.Range("I127").Value = 972.2
$ws.Range("K127").Value = 2916.6
$ws.Range("M127").Value = 2043.4
$ws.Range("H138").Value = 2510.359
$ws.Range("I138").Value = 1273.7273
$ws.Range("K138").Value = 3821.1819
$ws.Range("M138").Value = 1318.8181

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H31").Value = 9368.5
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H32").Value = 5115
$ws.Range("I32").Value = 2947.32
$ws.Range("K32").Value = 2947.32
$ws.Range("M32").Value = -2660.32
$ws.Range("H45").Value = 2057
$ws.Range("I45").Value = 1912.6923
$ws.Range("K45").Value = 1912.6923
$ws.Range("M45").Value = -1535.6923
$ws.Range("H97").Value = 696.6875
$ws.Range("I97").Value = 729.0769
$ws.Range("J97").Value = 556.3333
$ws.Range("K97").Value = 729.0769
$ws.Range("L97").Value = 556.3333
$ws.Range("M97").Value = -233.0769
$ws.Range("N97").Value = -1548.3333
$ws.Range("H102").Value = 3647951.8
$ws.Range("I102").Value = 3954261.8
$ws.Range("K102").Value = 3954261.8
$ws.Range("M102").Value = -3952639.8
$ws.Range("H122").Value = 25645026
$ws.Range("I122").Value = 55558892
$ws.Range("K122").Value = 166676676
$ws.Range("M122").Value = -166674226

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 46.666668
$ws.Range("I7").Value = 45.57143
$ws.Range("J7").Value = 50.5
$ws.Range("K7").Value = 45.57143
$ws.Range("L7").Value = 50.5
$ws.Range("M7").Value = 67.42857000000001
$ws.Range("N7").Value = -276.5
$ws.Range("H51").Value = 72333
$ws.Range("J51").Value = 80000
$ws.Range("L51").Value = 80000
$ws.Range("N51").Value = -81472
$ws.Range("H58").Value = 2296.6
$ws.Range("I58").Value = 2210.8572
$ws.Range("K58").Value = 2210.8572
$ws.Range("M58").Value = -2007.8572
$ws.Range("H61").Value = 72333
$ws.Range("J61").Value = 80000
$ws.Range("L61").Value = 80000
$ws.Range("N61").Value = -80696
$ws.Range("H62").Value = 6332.5
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
$ws.Range("H65").Value = 6332.5
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 100000
$ws.Range("N65").Value = -106240
$ws.Range("H134").Value = 2410.35
$ws.Range("I134").Value = 1598.1666
$ws.Range("J134").Value = 2758.4285
$ws.Range("K134").Value = 4794.4998
$ws.Range("L134").Value = 8275.2855
$ws.Range("M134").Value = -2259.4998
$ws.Range("N134").Value = -13345.2855
$ws.Range("H136").Value = 2296.6
$ws.Range("I136").Value = 2210.8572
$ws.Range("K136").Value = 6632.571599999999
$ws.Range("M136").Value = -4082.571599999999

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 3983.7144
$ws.Range("I80").Value = 3983.7144
$ws.Range("K80").Value = 3983.7144
$ws.Range("M80").Value = -2985.7144
$ws.Range("H83").Value = 3983.7144
$ws.Range("I83").Value = 3983.7144
$ws.Range("K83").Value = 19918.572
$ws.Range("M83").Value = -14926.572
$ws.Range("H122").Value = 2738.1875
$ws.Range("I122").Value = 2447.077
$ws.Range("K122").Value = 7341.231000000001
$ws.Range("M122").Value = -4891.231000000001

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 2190.32
$ws.Range("I46").Value = 1561.8
$ws.Range("K46").Value = 1561.8
$ws.Range("M46").Value = -1373.8
$ws.Range("H55").Value = 456.6316
$ws.Range("I55").Value = 289.08334
$ws.Range("J55").Value = 743.8570999999999
$ws.Range("K55").Value = 289.08334
$ws.Range("L55").Value = 743.8570999999999
$ws.Range("M55").Value = -116.08334
$ws.Range("N55").Value = -1089.8571
$ws.Range("H100").Value = 3099.2856
$ws.Range("I100").Value = 2898.75
$ws.Range("J100").Value = 3366.6667
$ws.Range("K100").Value = 2898.75
$ws.Range("L100").Value = 3366.6667
$ws.Range("M100").Value = -2357.75
$ws.Range("N100").Value = -4448.6667

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H62").Value = 2652280.5
$ws.Range("I62").Value = 5956257
$ws.Range("K62").Value = 5956257
$ws.Range("M62").Value = -5955633
$ws.Range("H65").Value = 2652280.5
$ws.Range("I65").Value = 5956257
$ws.Range("K65").Value = 29781285
$ws.Range("M65").Value = -29778165
$ws.Range("H70").Value = 37000
$ws.Range("I70").Value = 38000
$ws.Range("J70").Value = 36000
$ws.Range("K70").Value = 38000
$ws.Range("L70").Value = 36000
$ws.Range("M70").Value = -37685
$ws.Range("N70").Value = -36630
$ws.Range("H73").Value = 37000
$ws.Range("I73").Value = 38000
$ws.Range("J73").Value = 36000
$ws.Range("K73").Value = 38000
$ws.Range("L73").Value = 36000
$ws.Range("M73").Value = -36908
$ws.Range("N73").Value = -38184
$ws.Range("H81").Value = 9530399
$ws.Range("I81").Value = 3286
$ws.Range("K81").Value = 6572
$ws.Range("M81").Value = -5511
$ws.Range("H84").Value = 9530399
$ws.Range("I84").Value = 3286
$ws.Range("K84").Value = 32860
$ws.Range("M84").Value = -27556
$ws.Range("H96").Value = 2901.5
$ws.Range("J96").Value = 2949.5
$ws.Range("L96").Value = 2949.5
$ws.Range("N96").Value = -5695.5
$ws.Range("H122").Value = 1752
$ws.Range("I122").Value = 1541
$ws.Range("K122").Value = 4623
$ws.Range("M122").Value = -2173
$ws.Range("H126").Value = 1684.15
$ws.Range("I126").Value = 1678.5
$ws.Range("J126").Value = 1692.625
$ws.Range("K126").Value = 5035.5
$ws.Range("L126").Value = 5077.875
$ws.Range("M126").Value = -2565.5
$ws.Range("N126").Value = -10017.875
